# 🚌 141: 31/12 15:26 LP1912+6203+6173
# Append freshly scraped rows to the three schedule sheets and refresh the
# "last updated" / "total rows" header cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912": columns A..G = (blank), Hora_Scrap, Hora_Llegada,
#                 Línea, Minutos, Parada, Fecha
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 12:26:35"
$ws1.Range("A3").Value = "Total filas: 932"

$rows1 = @(
    @("12:26:24", "12:33", "23_HERNANDEZ", 7, "LP1912", "31/12/2025"),
    @("12:26:24", "12:40", "15X38_ABASTO", 14, "LP1912", "31/12/2025"),
    @("12:26:24", "12:51", "15_ABASTO", 25, "LP1912", "31/12/2025"),
    @("12:26:24", "12:54", "16_SANTA ANA", 28, "LP1912", "31/12/2025"),
    @("12:26:24", "12:58", "16_SANTA ANA", 32, "LP1912", "31/12/2025"),
    @("12:26:24", "13:01", "215C_EL PATO", 35, "LP1912", "31/12/2025"),
    @("12:26:24", "13:06", "14_ABASTO", 40, "LP1912", "31/12/2025"),
    @("12:26:24", "13:10", "23_HERNANDEZ", 44, "LP1912", "31/12/2025"),
    @("12:26:24", "13:11", "16_SANTA ANA", 45, "LP1912", "31/12/2025"),
    @("12:26:24", "13:21", "17_ROMERO", 55, "LP1912", "31/12/2025"),
    @("12:26:24", "13:30", "10_OLMOS", 64, "LP1912", "31/12/2025"),
    @("12:26:24", "13:31", "16_P MOR-SANTA ANA", 65, "LP1912", "31/12/2025"),
    @("12:26:24", "13:33", "23_HERNANDEZ", 67, "LP1912", "31/12/2025"),
    @("12:26:24", "13:51", "15_ABASTO", 85, "LP1912", "31/12/2025")
)

$startRow1 = 920
for ($i = 0; $i -lt $rows1.Length; $i++) {
    $row = $startRow1 + $i
    $data = $rows1[$i]
    $ws1.Cells.Item($row, 2).Value = $data[0]
    $ws1.Cells.Item($row, 3).Value = $data[1]
    $ws1.Cells.Item($row, 4).Value = $data[2]
    $ws1.Cells.Item($row, 5).Value = $data[3]
    $ws1.Cells.Item($row, 6).Value = $data[4]
    $ws1.Cells.Item($row, 7).Value = $data[5]
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": columns A..G = (blank), Fecha, Hora_Scrap,
#                      Hora_Llegada, Línea, Minutos, Parada
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 12:26:35"
$ws2.Range("A3").Value = "Total filas: 71"

$ws2.Cells.Item(72, 2).Value = "31/12/2025"
$ws2.Cells.Item(72, 3).Value = "12:26:24"
$ws2.Cells.Item(72, 4).Value = "13:01"
$ws2.Cells.Item(72, 5).Value = "215C_EL PATO"
$ws2.Cells.Item(72, 6).Value = 35
$ws2.Cells.Item(72, 7).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173": columns A..G = (blank), Fecha, Hora_Scrap,
#                     Hora_Llegada, Línea, Minutos, Parada
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 12:26:35"
$ws3.Range("A3").Value = "Total filas: 110"

$rows3 = @(
    @("31/12/2025", "12:26:34", "13:09", "215B_LP-P MOR-1 Y 57", 43, "L6173"),
    @("31/12/2025", "12:26:34", "13:14", "215A_LA PLATA", 48, "L6173"),
    @("31/12/2025", "12:26:29", "13:54", "215C_LA PLATA", 88, "L6203")
)

$startRow3 = 109
for ($i = 0; $i -lt $rows3.Length; $i++) {
    $row = $startRow3 + $i
    $data = $rows3[$i]
    $ws3.Cells.Item($row, 2).Value = $data[0]
    $ws3.Cells.Item($row, 3).Value = $data[1]
    $ws3.Cells.Item($row, 4).Value = $data[2]
    $ws3.Cells.Item($row, 5).Value = $data[3]
    $ws3.Cells.Item($row, 6).Value = $data[4]
    $ws3.Cells.Item($row, 7).Value = $data[5]
}
